$d = $word.ActiveDocument

# Paragraph 2 ("The THEN paragraph.") keeps its <w:pPr> (the tab stop),
# but all of its runs (and the spell-check proofErr markers around them)
# are replaced by a broken/invalid "m:if" field-code run sequence
# followed by a bold red error-message run. Paragraphs 3 ("End of
# demonstration.") and 4 (empty) are dropped entirely.

$p2 = $d.Paragraphs.Item(2)
$p4 = $d.Paragraphs.Item(4)

# A range spanning from the start of paragraph 2's content through the
# end of paragraph 4 (its paragraph mark included) crosses three
# paragraph marks (end of p2, end of p3, end of p4). Replacing that whole
# span via InsertXML collapses paragraphs 2-4 into the single new
# paragraph described below (so we must restate p2's own <w:pPr> inside
# the replacement markup).
$full = $d.Range($p2.Range.Start, $p4.Range.End)

$runs = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText xml:space="preserve">m:if </w:instrText></w:r>' +
        '<w:r><w:instrText xml:space="preserve">self.name </w:instrText></w:r>' +
        '<w:r><w:instrText>=</w:instrText></w:r>' +
        '<w:r><w:instrText xml:space="preserve"> </w:instrText></w:r>' +
        '<w:r><w:instrText>''</w:instrText></w:r>' +
        '<w:r><w:instrText>anydsl</w:instrText></w:r>' +
        '<w:r><w:instrText>''</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '<w:r><w:rPr><w:b w:val="true"/><w:color w:val="FF0000"/></w:rPr>' +
        '<w:t>Invalid if statement: m:elseif, m:else or m:endif expected here.</w:t></w:r>'

$paragraph = '<w:p><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr>' + $runs + '</w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $paragraph + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($xml)
